$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the values in row 3 (C3:J3) while keeping their formatting/style
$ws.Range("C3:J3").ClearContents()

# Update G30 value from 0 to 5 (dependent formula in L30 recalculates automatically)
$ws.Range("G30").Value = 5

# Move the active selection to C4 (was O14)
$ws.Range("C4").Select()
